$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author removed the "grouping1" / "grouping2" header columns from the
# experimental-design template (commit: "Added group column to experimental
# file: grouping1 and grouping2" -- i.e. those two columns, which lived at
# AA:AB, are the ones being dropped here, shifting every later column left
# by two positions).
#
# Mimic the manual workflow: select the target columns, then delete them so
# the remaining data/headers shift left.
$ws.Range("T1").Select() | Out-Null
$ws.Range("AA1:AB1048576").Select() | Out-Null
$ws.Range("AA:AB").Delete()
